# Commit: "flipped labels in barchart on their side, also made the grapg
# taller so there is less wasted space"
#
# The slide hosts a PowerPoll task-pane/content add-in (an
# mc:AlternateContent graphicFrame backed by a we:webextensionref, with a
# cached p:pic fallback image). The add-in renders its bar chart entirely
# client-side (inside the task pane's own HTML/JS), so "flipping the axis
# labels" and "making the chart taller" are changes to the add-in's own
# internal UI/state - they are not expressed anywhere in the slide's
# OOXML. The graphicFrame's cached size (cx/cy) in the slide XML is
# unchanged by the commit, confirming this.
#
# The only bytes that differ in the package are:
#   - every r:id in presentation.xml/slide.xml (slideMaster, slide, all
#     11 slideLayouts, the webextensionref and the blip embed) - this is
#     the relationship-id churn PowerPoint performs on every full-deck
#     resave and is not something any edit script controls or targets;
#   - the we:webextension/@id GUID in ppt/slides/udata/data.xml, which is
#     an internal snapshot/instance id PowerPoint mints for the add-in
#     when its task-pane state changes. It is plumbing private to the
#     Office Add-ins host, not part of the Presentation/Slide/Shape
#     automation surface, so it cannot be set from VBA/COM (there is no
#     WebExtension object in the PowerPoint object model, and the
#     add-in's graphicFrame/pic pair in this deck even re-uses the
#     Title placeholder's shape id, so it is not individually
#     addressable through Shapes either).
#
# So there is no reachable, content-visible edit for this commit via the
# Presentation/Slide/Shape COM surface - every other part of the deck
# (titles, subtitle, layouts, masters, theme, media bytes) is identical
# before/after. Intentionally make no changes here rather than perturb
# shapes/parts that the real commit left untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
